$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain stored as text,
# matching the original inline-string cell type. We temporarily switch
# their number format to Text ("@") before assigning the value, then
# restore the default "Normal" style so the saved XML has no explicit
# style index (matching the untouched cells around them).
$textCells = @(
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D20",
    "D22",
    "D23",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)

foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "27.379.64"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "1.832.36"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  -1.08%  "
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D7").Value = "0.4257"
$ws.Range("E7").Value = "  -1.51%  "
$ws.Range("D8").Value = "0.3705"
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("D9").Value = "0.07266"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").Value = "0.8680"
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("D11").Value = "21.17"
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("D12").Value = "1.832.04"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("D13").Value = "6.742"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "0.07115"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "5.323"
$ws.Range("E15").Value = "  -3.16%  "
$ws.Range("D16").Value = "89.15"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("D18").Value = "0.000008876"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").Value = "15.12"
$ws.Range("D21").Value = "27.416.28"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("D22").Value = "5.144"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").Value = "10.92"
$ws.Range("E23").Value = "  -2.57%  "
$ws.Range("D24").Value = "2.048.14"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").Value = "2.003"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").Value = "152.85"
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("D27").Value = "2.201"
$ws.Range("E27").Value = "  +4.56%  "
$ws.Range("D28").Value = "18.48"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D29").Value = "5.259"
$ws.Range("E29").Value = "  -3.23%  "
$ws.Range("D30").Value = "116.70"
$ws.Range("E30").Value = "  -4.40%  "
$ws.Range("D31").Value = "0.08894"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").Value = "1.208"
$ws.Range("E32").Value = "  -2.85%  "
$ws.Range("D33").Value = "0.7631"
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("D34").Value = "4.494"
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").Value = "2.818"
$ws.Range("E35").Value = "  -3.80%  "
$ws.Range("D36").Value = "1.004"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").Value = "1.124"
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("D38").Value = "0.01984"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "0.05301"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "7.262"
$ws.Range("E40").Value = "  +2.55%  "
$ws.Range("D41").Value = "2.890"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").Value = "0.1710"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").Value = "0.5099"
$ws.Range("D44").Value = "8.742"
$ws.Range("E44").Value = "  -1.58%  "
$ws.Range("D45").Value = "10.71"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "108.14"
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("D47").Value = "0.4790"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").Value = "0.06398"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("D50").Value = "1.675"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").Value = "1.861"
$ws.Range("E51").Value = "  -2.15%  "

# Restore default style on the cells we temporarily reformatted
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
